$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking"): Right count 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total 135 -> 108, Wrong total -1 -> -2, Max text updated
$ws.Range("B12").Value = 108
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "106 / 112"
